# Auto-generated edit script for Top50_DataComp.xlsx
# Applies the monthly data refresh: M2/FX length counters bumped by one
# and M2/FX first/last-observation dates rolled forward, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to build a clean midnight DateTime so Excel stores an integer serial
# (no time-of-day fraction) for each date cell.
function New-CleanDate([int]$y, [int]$m, [int]$d) {
    return Get-Date -Year $y -Month $m -Day $d -Hour 0 -Minute 0 -Second 0 -Millisecond 0
}

# --- Integer counter updates (M2_Len / FX_Len columns C & D) ---
$ws.Range("C6").Value = 443
$ws.Range("D8").Value = 410
$ws.Range("D10").Value = 483
$ws.Range("C12").Value = 371
$ws.Range("D12").Value = 352
$ws.Range("C13").Value = 467
$ws.Range("D14").Value = 396
$ws.Range("C16").Value = 455
$ws.Range("D16").Value = 410
$ws.Range("C17").Value = 370
$ws.Range("D17").Value = 394
$ws.Range("D19").Value = 398
$ws.Range("C20").Value = 481
$ws.Range("C21").Value = 310
$ws.Range("C22").Value = 323
$ws.Range("D22").Value = 367
$ws.Range("D23").Value = 402
$ws.Range("C24").Value = 216
$ws.Range("D24").Value = 410
$ws.Range("D25").Value = 315
$ws.Range("C26").Value = 333
$ws.Range("D26").Value = 313
$ws.Range("D28").Value = 380
$ws.Range("C29").Value = 262
$ws.Range("D29").Value = 235
$ws.Range("D30").Value = 217
$ws.Range("C31").Value = 394
$ws.Range("D33").Value = 398
$ws.Range("C34").Value = 202
$ws.Range("D34").Value = 317
$ws.Range("C35").Value = 404
$ws.Range("D35").Value = 317
$ws.Range("D36").Value = 410
$ws.Range("C37").Value = 466
$ws.Range("D37").Value = 317
$ws.Range("C38").Value = 359
$ws.Range("D38").Value = 367
$ws.Range("C39").Value = 227
$ws.Range("D39").Value = 223
$ws.Range("D40").Value = 314
$ws.Range("D41").Value = 316
$ws.Range("D42").Value = 218
$ws.Range("D43").Value = 317
$ws.Range("D44").Value = 304
$ws.Range("D45").Value = 317
$ws.Range("C46").Value = 334
$ws.Range("D46").Value = 298
$ws.Range("C47").Value = 335
$ws.Range("D47").Value = 260
$ws.Range("D48").Value = 315
$ws.Range("D49").Value = 313
$ws.Range("D50").Value = 238
$ws.Range("D51").Value = 317
$ws.Range("D52").Value = 315

# --- Date rollforwards (M2_1stDate/M2_LastDate/FX_1stDate/FX_LastDate columns E-H) ---
$ws.Range("G2").Value = New-CleanDate 1982 5 3
$ws.Range("H2").Value = New-CleanDate 2023 12 1
$ws.Range("E3").Value = New-CleanDate 1982 3 1
$ws.Range("F3").Value = New-CleanDate 2023 10 1
$ws.Range("E4").Value = New-CleanDate 1982 3 1
$ws.Range("F4").Value = New-CleanDate 2023 10 1
$ws.Range("G4").Value = New-CleanDate 1982 5 3
$ws.Range("H4").Value = New-CleanDate 2023 12 1
$ws.Range("G5").Value = New-CleanDate 1982 5 3
$ws.Range("H5").Value = New-CleanDate 2023 12 1
$ws.Range("F6").Value = New-CleanDate 2023 10 1
$ws.Range("G6").Value = New-CleanDate 1982 5 3
$ws.Range("H6").Value = New-CleanDate 2023 12 1
$ws.Range("G7").Value = New-CleanDate 1982 5 3
$ws.Range("H7").Value = New-CleanDate 2023 12 1
$ws.Range("E8").Value = New-CleanDate 1982 3 1
$ws.Range("F8").Value = New-CleanDate 2023 10 1
$ws.Range("H8").Value = New-CleanDate 2023 12 1
$ws.Range("E9").Value = New-CleanDate 1982 3 1
$ws.Range("F9").Value = New-CleanDate 2023 10 1
$ws.Range("G9").Value = New-CleanDate 1982 5 3
$ws.Range("H9").Value = New-CleanDate 2023 12 1
$ws.Range("E10").Value = New-CleanDate 1982 3 1
$ws.Range("F10").Value = New-CleanDate 2023 10 1
$ws.Range("H10").Value = New-CleanDate 2023 12 1
$ws.Range("G11").Value = New-CleanDate 1982 5 3
$ws.Range("H11").Value = New-CleanDate 2023 12 1
$ws.Range("F12").Value = New-CleanDate 2023 10 1
$ws.Range("H12").Value = New-CleanDate 2023 12 1
$ws.Range("F13").Value = New-CleanDate 2023 10 1
$ws.Range("G13").Value = New-CleanDate 1982 5 3
$ws.Range("H13").Value = New-CleanDate 2023 12 1
$ws.Range("H14").Value = New-CleanDate 2023 12 1
$ws.Range("G15").Value = New-CleanDate 1982 4 1
$ws.Range("H15").Value = New-CleanDate 2023 12 1
$ws.Range("F16").Value = New-CleanDate 2023 10 1
$ws.Range("H16").Value = New-CleanDate 2023 12 1
$ws.Range("F17").Value = New-CleanDate 2023 10 1
$ws.Range("H17").Value = New-CleanDate 2023 12 1
$ws.Range("E18").Value = New-CleanDate 1982 3 1
$ws.Range("F18").Value = New-CleanDate 2023 10 1
$ws.Range("G18").Value = New-CleanDate 1982 5 3
$ws.Range("H18").Value = New-CleanDate 2023 12 1
$ws.Range("E19").Value = New-CleanDate 1982 3 1
$ws.Range("F19").Value = New-CleanDate 2023 10 1
$ws.Range("H19").Value = New-CleanDate 2023 12 1
$ws.Range("F20").Value = New-CleanDate 2023 10 1
$ws.Range("G20").Value = New-CleanDate 1982 5 3
$ws.Range("H20").Value = New-CleanDate 2023 12 1
$ws.Range("F21").Value = New-CleanDate 2023 10 1
$ws.Range("G21").Value = New-CleanDate 1982 5 3
$ws.Range("H21").Value = New-CleanDate 2023 12 1
$ws.Range("F22").Value = New-CleanDate 2023 10 1
$ws.Range("H22").Value = New-CleanDate 2023 12 1
$ws.Range("H23").Value = New-CleanDate 2023 12 1
$ws.Range("F24").Value = New-CleanDate 2023 11 1
$ws.Range("H24").Value = New-CleanDate 2023 12 1
$ws.Range("H25").Value = New-CleanDate 2023 12 1
$ws.Range("F26").Value = New-CleanDate 2023 10 1
$ws.Range("H26").Value = New-CleanDate 2023 12 1
$ws.Range("E27").Value = New-CleanDate 1982 3 1
$ws.Range("F27").Value = New-CleanDate 2023 10 1
$ws.Range("G27").Value = New-CleanDate 1982 5 3
$ws.Range("H27").Value = New-CleanDate 2023 12 1
$ws.Range("H28").Value = New-CleanDate 2023 12 1
$ws.Range("F29").Value = New-CleanDate 2023 10 1
$ws.Range("H29").Value = New-CleanDate 2023 12 1
$ws.Range("E30").Value = New-CleanDate 1982 3 1
$ws.Range("F30").Value = New-CleanDate 2023 10 1
$ws.Range("H30").Value = New-CleanDate 2023 12 1
$ws.Range("F31").Value = New-CleanDate 2023 10 1
$ws.Range("G31").Value = New-CleanDate 1982 5 3
$ws.Range("H31").Value = New-CleanDate 2023 12 1
$ws.Range("E32").Value = New-CleanDate 1982 3 1
$ws.Range("F32").Value = New-CleanDate 2023 10 1
$ws.Range("G32").Value = New-CleanDate 1982 5 3
$ws.Range("H32").Value = New-CleanDate 2023 12 1
$ws.Range("H33").Value = New-CleanDate 2023 12 1
$ws.Range("F34").Value = New-CleanDate 2023 10 1
$ws.Range("H34").Value = New-CleanDate 2023 12 1
$ws.Range("F35").Value = New-CleanDate 2023 9 1
$ws.Range("H35").Value = New-CleanDate 2023 12 1
$ws.Range("H36").Value = New-CleanDate 2023 12 1
$ws.Range("F37").Value = New-CleanDate 2023 10 1
$ws.Range("H37").Value = New-CleanDate 2023 12 1
$ws.Range("F38").Value = New-CleanDate 2023 10 1
$ws.Range("H38").Value = New-CleanDate 2023 12 1
$ws.Range("F39").Value = New-CleanDate 2023 10 1
$ws.Range("H39").Value = New-CleanDate 2023 12 1
$ws.Range("H40").Value = New-CleanDate 2023 12 1
$ws.Range("H41").Value = New-CleanDate 2023 12 1
$ws.Range("H42").Value = New-CleanDate 2023 12 1
$ws.Range("H43").Value = New-CleanDate 2023 12 1
$ws.Range("H44").Value = New-CleanDate 2023 12 1
$ws.Range("H45").Value = New-CleanDate 2023 12 1
$ws.Range("F46").Value = New-CleanDate 2023 10 1
$ws.Range("H46").Value = New-CleanDate 2023 12 1
$ws.Range("F47").Value = New-CleanDate 2023 10 1
$ws.Range("H47").Value = New-CleanDate 2023 12 1
$ws.Range("H48").Value = New-CleanDate 2023 12 1
$ws.Range("H49").Value = New-CleanDate 2023 12 1
$ws.Range("H50").Value = New-CleanDate 2023 12 1
$ws.Range("E51").Value = New-CleanDate 1981 9 1
$ws.Range("F51").Value = New-CleanDate 2023 9 1
$ws.Range("H51").Value = New-CleanDate 2023 12 1
$ws.Range("H52").Value = New-CleanDate 2023 12 1
